$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2e5d215f-0658-4e53-afe1-d7b6a35800a7"
$ws.Range("D3").Value = "dc12476e-0041-4bf3-b221-2945de51ae3d"
$ws.Range("D4").Value = "8b521257-4f8e-4660-989c-858519fab6a1"

$ws.Range("C2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6ImRhbmllbDVmIiwicGFzc3dvcmQiOiJBejI1Mjg4QCIsImlhdCI6MTcwMjgyMzMzNX0.SXAT0MmUMP0xjskTleFUU5bScEQmpuwTMLDhcWTMcIg"
$ws.Range("C3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6IkpvcmdlMjUyNSIsInBhc3N3b3JkIjoiYXNUMzU2NDQ0QCIsImlhdCI6MTcwMjgyMzMzNn0._MOQkwJA9OcY0W63cJqWmSilvzH6qvXvGj6YgeVr9ak"
$ws.Range("C4").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6Im1hcmlvMzUiLCJwYXNzd29yZCI6Im1BcmlvdXVnQDMiLCJpYXQiOjE3MDI4MjMzMzd9.xjEQgL5uO_0ke5Dkybz0QD2IYchEXW9JbrJ8geGblgI"
